$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 90bdfd9a... row
$wsOverview.Range("G4").Value = "2016-10-27 09:20:17"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 90bdfd9a... row
$wsZhCn.Range("H4").Value = "2016-10-27 09:20:04"
$wsZhCn.Range("K4").Value = "2016-10-27 09:20:46"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 90bdfd9a... row
$wsDeDe.Range("H4").Value = "2016-10-27 09:20:17"
$wsDeDe.Range("K4").Value = "2016-10-27 09:21:04"
